# Increase MaxInvest Storage Adapt Szenarios Existing Units
#
# On the "Power Storage" sheet, column S holds "MaxInvest" for the
# existing-unit rows (7-11). Raise the cap from 8 to 15 for each of them,
# then leave the selection on that range to mirror the saved workbook
# state (S8:S11 active on the frozen/bottom-left pane).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S7").Value = 15
$ws.Range("S8").Value = 15
$ws.Range("S9").Value = 15
$ws.Range("S10").Value = 15
$ws.Range("S11").Value = 15

$ws.Range("S8:S11").Select()
